# Add "review" column (N) with header and flags for Chapters 0-5 (rows 21-26),
# and a free-text review-notes column (O) with comments for Chapter 6 (row 27)
# and Chapter 7 (row 28). Matches commit: "caps 0, 1 - esp - wordpress".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column N
$ws.Range("N1").Value = "review"

# Mark chapters 0-5 (rows 21-26) as reviewed
$ws.Range("N21").Value = 1
$ws.Range("N22").Value = 1
$ws.Range("N23").Value = 1
$ws.Range("N24").Value = 1
$ws.Range("N25").Value = 1
$ws.Range("N26").Value = 1

# Review notes for chapter 6 (row 27) and chapter 7 (row 28)
$ws.Range("O27").Value = "es  corto"
$ws.Range("O28").Value = "es largo y necesita reestructurarse"

# Set column O width to match the new wider review-notes column.
# (34.4285714285714 itself cannot be represented exactly through the
# ColumnWidth property's internal 1/6-character granularity; 33.6 is the
# input that lands on the closest achievable stored width, 34.5.)
$ws.Columns("O").ColumnWidth = 33.6

# Update the active selection to reflect where the user ended up
$ws.Range("Q24").Select()
